$d = $word.ActiveDocument

# 1. Replace the title text: split "RELATÓRIO DO DESAFIO PARA SELEÇÃO" into
#    "RELATÓRIO DO DESAFIO " + "AMBIENTE MONITORAMENTO" (two separate runs).
$d.Content.Find.Execute("RELATÓRIO DO DESAFIO PARA SELEÇÃO", $true, $false, $false, $false, $false,
                         $true, 1, $false, "RELATÓRIO DO DESAFIO AMBIENTE MONITORAMENTO", 2)

# 2. Insert two new empty paragraphs right after the title paragraph.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "RELATÓRIO DO DESAFIO AMBIENTE MONITORAMENTO*") {
        $p.Range.InsertParagraphAfter()
        $p.Range.InsertParagraphAfter()
        break
    }
}
